$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 270, shifting rows 270:330 down to 271:331
$ws.Rows("270:270").Insert()

# Populate the newly inserted row 270 with the new weekly record
$ws.Range("A270").Value2 = 9
$ws.Range("B270").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C270").Value2 = "Metropolitana"
$ws.Range("D270").Value2 = 45204
$ws.Range("E270").Value2 = 13
$ws.Range("F270").Value2 = "Fruta"
$ws.Range("G270").Value2 = 100101
$ws.Range("H270").Value2 = "Berries"
$ws.Range("I270").Value2 = 100101001
$ws.Range("J270").Value2 = "Arándano (blue)"
$ws.Range("K270").Value2 = "Sin especificar"
$ws.Range("L270").Value2 = "Primera"
$ws.Range("M270").Value2 = 52
$ws.Range("N270").Value2 = 11000
$ws.Range("O270").Value2 = 12000
$ws.Range("P270").Value2 = 11500
$ws.Range("Q270").Value2 = "`$/bandeja 2 kilos"
$ws.Range("R270").Value2 = "Provincia de Limarí"
$ws.Range("S270").Value2 = 5750
$ws.Range("T270").Value2 = 2
